# Agregando riesgo identificado en proceso
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New risk row (row 14): "Acceso denegado al dispositivo de respaldo"
$ws.Range("C14").Value = "Acceso denegado al dispositivo de respaldo"
$ws.Range("D14").Value = 4
$ws.Range("E14").Value = 0.15
$ws.Range("G14").Value = 4
$ws.Range("H14").Value = "Generar adquisicion de un dispositivo secundario que lea discos duros del actual"
$ws.Range("I14").Value = "Solicitar al usuario cargar en usb la informacion mas reciente que ha generado"
$ws.Range("J14").Value = "Jovanny Zepeda"
$ws.Range("K14").Value = "Ocurrido"

$ws.Rows.Item(14).RowHeight = 45

# Update view / selection to match the saved workbook state
$ws.Range("K14").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 4
$win.Zoom = 85
